$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "만남" intent now also carries entity info and a templated query.
$ws.Range("B3").Value = "B_DT,B_TI"
$ws.Range("C3").Value = "{B_DT} 만남"
$ws.Range("D3").WrapText = $true

# Row 5 used to hold "주문" - it now holds "욕설" (profanity) together with
# its warning answer that used to live further down in row 7.
$ws.Range("A5").Value = "욕설"
$ws.Range("D5").Value = "[경고]" + [char]10 + "상대방과 나를 위해 욕설 사용은 자제해주세요."
$ws.Range("D5").WrapText = $true
$ws.Rows("5:5").RowHeight = 51.75

# Row 6 ("예약") and row 7 ("욕설" + its warning, now moved up to row 5) are cleared.
$ws.Range("A6").ClearContents()
$ws.Range("A7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Rows("7:7").RowHeight = 17.25

$ws.Range("E3").Select()
